$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.072147192379749
$ws.Range("D2").Value = 1.076565569718524
$ws.Range("E2").Value = 1.085088929233691
$ws.Range("F2").Value = 1.090639572145684
$ws.Range("I2").Value = 1.058197580245922
$ws.Range("J2").Value = 1.077068206821798
$ws.Range("K2").Value = 1.079249166068432
$ws.Range("L2").Value = 1.08775026307151
$ws.Range("M2").Value = 1.09328660986049
$ws.Range("N2").Value = 1.029217183606594
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.073276496742265
$ws.Range("D3").Value = 1.077484330725824
$ws.Range("E3").Value = 1.086166185623608
$ws.Range("F3").Value = 1.091702139955989
$ws.Range("I3").Value = 1.05855015657959
$ws.Range("J3").Value = 1.077854854150782
$ws.Range("K3").Value = 1.079985207602885
$ws.Range("L3").Value = 1.088646024096477
$ws.Range("M3").Value = 1.09416875099252
$ws.Range("N3").Value = 1.029488737811342
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.074007473202132
$ws.Range("D4").Value = 1.078079010712378
$ws.Range("E4").Value = 1.086863782987127
$ws.Range("F4").Value = 1.092390201083886
$ws.Range("I4").Value = 1.058777248032534
$ws.Range("J4").Value = 1.078363502905558
$ws.Range("K4").Value = 1.080461021134585
$ws.Range("L4").Value = 1.089225576487071
$ws.Range("M4").Value = 1.094739454712172
$ws.Range("N4").Value = 1.029664134866417
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.074314833834728
$ws.Range("D5").Value = 1.078329057073278
$ws.Range("E5").Value = 1.087157182227498
$ws.Range("F5").Value = 1.092679583630534
$ws.Range("I5").Value = 1.058872466109228
$ws.Range("J5").Value = 1.078577251549527
$ws.Range("K5").Value = 1.080660944249173
$ws.Range("L5").Value = 1.089469204744816
$ws.Range("M5").Value = 1.094979354403455
$ws.Range("N5").Value = 1.029737795839479
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.074366444458669
$ws.Range("D6").Value = 1.078371043465952
$ws.Range("E6").Value = 1.087206452831616
$ws.Range("F6").Value = 1.092728179363118
$ws.Range("I6").Value = 1.058888438913601
$ws.Range("J6").Value = 1.078613135768881
$ws.Range("K6").Value = 1.080694505832195
$ws.Range("L6").Value = 1.08951011008758
$ws.Range("M6").Value = 1.095019633198866
$ws.Range("N6").Value = 1.029750159383343
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.074011579940988
$ws.Range("D7").Value = 1.078082351675837
$ws.Range("E7").Value = 1.086867702894124
$ws.Range("F7").Value = 1.092394067349801
$ws.Range("I7").Value = 1.058778521328098
$ws.Range("J7").Value = 1.078366359368379
$ws.Range("K7").Value = 1.080463692944988
$ws.Range("L7").Value = 1.089228831920974
$ws.Range("M7").Value = 1.094742660358574
$ws.Range("N7").Value = 1.02966511942654
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.072528796426957
$ws.Range("D8").Value = 1.076876031850346
$ws.Range("E8").Value = 1.085452881126007
$ws.Range("F8").Value = 1.090998566499567
$ws.Range("I8").Value = 1.058316952251816
$ws.Range("J8").Value = 1.077334133482134
$ws.Range("K8").Value = 1.079498008826855
$ws.Range("L8").Value = 1.088053003119794
$ws.Range("M8").Value = 1.09358475436727
$ws.Range("N8").Value = 1.029309022108801
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.06991776900602
$ws.Range("D9").Value = 1.074751732411479
$ws.Range("E9").Value = 1.082963926468196
$ws.Range("F9").Value = 1.088543417366475
$ws.Range("I9").Value = 1.057495575778267
$ws.Range("J9").Value = 1.075512425878499
$ws.Range("K9").Value = 1.077792871505257
$ws.Range("L9").Value = 1.085980546617515
$ws.Range("M9").Value = 1.091543608630348
$ws.Range("N9").Value = 1.028679113689714
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.068178282398255
$ws.Range("D10").Value = 1.073336478880764
$ws.Range("E10").Value = 1.081307409247907
$ws.Range("F10").Value = 1.08690927927142
$ws.Range("I10").Value = 1.056942591208317
$ws.Range("J10").Value = 1.074296072288738
$ws.Range("K10").Value = 1.07665377610398
$ws.Range("L10").Value = 1.084598574307084
$ws.Range("M10").Value = 1.090182330138374
$ws.Range("N10").Value = 1.028257555166844
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.067425342250602
$ws.Range("D11").Value = 1.072723884478004
$ws.Range("E11").Value = 1.080590779034445
$ws.Range("F11").Value = 1.086202302835113
$ws.Range("I11").Value = 1.056701861578635
$ws.Range("J11").Value = 1.07376892984533
$ws.Range("K11").Value = 1.076159980449362
$ws.Range("L11").Value = 1.08400008454144
$ws.Range("M11").Value = 1.089592758880156
$ws.Range("N11").Value = 1.028074632759425
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.067145706425516
$ws.Range("D12").Value = 1.072496372682541
$ws.Range("E12").Value = 1.08032468854973
$ws.Range("F12").Value = 1.085939792773892
$ws.Range("I12").Value = 1.056612250924414
$ws.Range("J12").Value = 1.073573057334769
$ws.Range("K12").Value = 1.075976478516746
$ws.Range("L12").Value = 1.083777765579345
$ws.Range("M12").Value = 1.089373746377987
$ws.Range("N12").Value = 1.028006629485345
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.067205687476048
$ws.Range("D13").Value = 1.072545173239257
$ws.Range("E13").Value = 1.080381761435232
$ws.Range("F13").Value = 1.085996097890519
$ws.Range("I13").Value = 1.056631481446894
$ws.Range("J13").Value = 1.073615075749223
$ws.Range("K13").Value = 1.076015844118353
$ws.Range("L13").Value = 1.083825454352594
$ws.Range("M13").Value = 1.089420726187371
$ws.Range("N13").Value = 1.028021219030666
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.067402226661114
$ws.Range("D14").Value = 1.072705077607792
$ws.Range("E14").Value = 1.080568781916144
$ws.Range("F14").Value = 1.086180601785884
$ws.Range("I14").Value = 1.056694458270692
$ws.Range("J14").Value = 1.073752740357864
$ws.Range("K14").Value = 1.076144813841599
$ws.Range("L14").Value = 1.083981707857412
$ws.Range("M14").Value = 1.089574655636788
$ws.Range("N14").Value = 1.028069012765242
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.067523326174058
$ws.Range("D15").Value = 1.072803604346053
$ws.Range("E15").Value = 1.08068402437243
$ws.Range("F15").Value = 1.086294292974732
$ws.Range("I15").Value = 1.056733234799607
$ws.Range("J15").Value = 1.073837551004705
$ws.Range("K15").Value = 1.07622426517807
$ws.Range("L15").Value = 1.08407797904167
$ws.Range("M15").Value = 1.089669494064271
$ws.Range("N15").Value = 1.028098452412291
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.068228258160294
$ws.Range("D16").Value = 1.073377139385571
$ws.Range("E16").Value = 1.081354983371335
$ws.Range("F16").Value = 1.08695621194169
$ws.Range("I16").Value = 1.056958540587513
$ws.Range("J16").Value = 1.074331047426054
$ws.Range("K16").Value = 1.076686535889001
$ws.Range("L16").Value = 1.084638292259362
$ws.Range("M16").Value = 1.090221455316709
$ws.Range("N16").Value = 1.028269687022916
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.068670514869079
$ws.Range("D17").Value = 1.073736961609462
$ws.Range("E17").Value = 1.081776032842032
$ws.Range("F17").Value = 1.087371581107581
$ws.Range("I17").Value = 1.057099525214808
$ws.Range("J17").Value = 1.074640483164839
$ws.Range("K17").Value = 1.076976356116873
$ws.Range("L17").Value = 1.084989738601374
$ws.Range("M17").Value = 1.09056765146619
$ws.Range("N17").Value = 1.028376994996275
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.068928501732629
$ws.Range("D18").Value = 1.073946861022113
$ws.Range("E18").Value = 1.082021686859765
$ws.Range("F18").Value = 1.087613918370201
$ws.Range("I18").Value = 1.057181635314161
$ws.Range("J18").Value = 1.074820928073331
$ws.Range("K18").Value = 1.077145349345801
$ws.Range("L18").Value = 1.085194722960841
$ws.Range("M18").Value = 1.090769569416569
$ws.Range("N18").Value = 1.028439548791268
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.069016473020317
$ws.Range("D19").Value = 1.074018434891994
$ws.Range("E19").Value = 1.082105459220176
$ws.Range("F19").Value = 1.087696559202068
$ws.Range("I19").Value = 1.057209611732223
$ws.Range("J19").Value = 1.074882447714948
$ws.Range("K19").Value = 1.077202962513284
$ws.Range("L19").Value = 1.085264615920649
$ws.Range("M19").Value = 1.090838416119937
$ws.Range("N19").Value = 1.028460871718977
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.068623062224859
$ws.Range("D20").Value = 1.07369835389731
$ws.Range("E20").Value = 1.081730851704884
$ws.Range("F20").Value = 1.087327009777948
$ws.Range("I20").Value = 1.05708441170943
$ws.Range("J20").Value = 1.074607288162548
$ws.Range("K20").Value = 1.076945266736273
$ws.Range("L20").Value = 1.08495203259951
$ws.Range("M20").Value = 1.090530509179968
$ws.Range("N20").Value = 1.02836548570637
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.067344349680615
$ws.Range("D21").Value = 1.072657988844403
$ws.Range("E21").Value = 1.080513706345161
$ws.Range("F21").Value = 1.08612626743566
$ws.Range("I21").Value = 1.056675918492894
$ws.Range("J21").Value = 1.073712203450073
$ws.Range("K21").Value = 1.076106837799084
$ws.Range("L21").Value = 1.083935695461159
$ws.Range("M21").Value = 1.089529327789653
$ws.Range("N21").Value = 1.028054940295247
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.066540601505886
$ws.Range("D22").Value = 1.072004061228333
$ws.Range("E22").Value = 1.079749004995344
$ws.Range("F22").Value = 1.085371847994411
$ws.Range("I22").Value = 1.056417966192816
$ws.Range("J22").Value = 1.073149032681885
$ws.Range("K22").Value = 1.07557919679745
$ws.Range("L22").Value = 1.083296607407208
$ws.Range("M22").Value = 1.08889973249816
$ws.Range("N22").Value = 1.027859353992392
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.06696666214838
$ws.Range("D23").Value = 1.072350702497308
$ws.Range("E23").Value = 1.080154334032884
$ws.Range("F23").Value = 1.08577172934332
$ws.Range("I23").Value = 1.056554817440911
$ws.Range("J23").Value = 1.07344761772538
$ws.Range("K23").Value = 1.075858955524065
$ws.Range("L23").Value = 1.083635407399164
$ws.Range("M23").Value = 1.089233503618261
$ws.Range("N23").Value = 1.027963069638007
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.068644503951574
$ws.Range("D24").Value = 1.073715798994972
$ws.Range("E24").Value = 1.081751266921187
$ws.Range("F24").Value = 1.087347149459402
$ws.Range("I24").Value = 1.057091241234544
$ws.Range("J24").Value = 1.074622287690162
$ws.Range("K24").Value = 1.076959314855878
$ws.Range("L24").Value = 1.084969070344346
$ws.Range("M24").Value = 1.090547292218527
$ws.Range("N24").Value = 1.028370686373761
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.070592568801119
$ws.Range("D25").Value = 1.075300748816841
$ws.Range("E25").Value = 1.08360688995554
$ws.Range("F25").Value = 1.089177669417959
$ws.Range("I25").Value = 1.057708873164795
$ws.Range("J25").Value = 1.075983712399563
$ws.Range("K25").Value = 1.078234102442033
$ws.Range("L25").Value = 1.086516385201526
$ws.Range("M25").Value = 1.092071384546279
$ws.Range("N25").Value = 1.028842246097705
